$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new Information_Ratio test rows (rows 62-63) following the
# existing Test / Description / macro column pattern.
$ws.Range("A62").Value = "Information_Ratio1"
$ws.Range("B62").Value = "Test Information Ratio with scale=252"
$ws.Range("C62").Value = "Information_Ratio_test1"

$ws.Range("A63").Value = "Information_Ratio2"
$ws.Range("B63").Value = "Test Information Ratio with scale=1"
$ws.Range("C63").Value = "Information_Ratio_test2"

# Update the active selection to match the new last-used cell.
$ws.Range("C63").Select()
